$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXBuySell")
Write-Output $ws.PageSetup.Orientation
